$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide the CitationGeoKeys requirement rows (2-6) ---
$ws.Rows.Item(2).Hidden = $false
$ws.Rows.Item(3).Hidden = $false
$ws.Rows.Item(4).Hidden = $false
$ws.Rows.Item(5).Hidden = $false
$ws.Rows.Item(6).Hidden = $false

# --- Status updates: mark the AsciiParams requirement rows (3 & 4) as Done ---
# (row 2, 5, 6 keep their existing / still-empty Status column)
$ws.Range("D3").Value = "D"
$ws.Range("D4").Value = "D"

# --- Hide the Ellipsoid* requirement rows (7-20), now filtered out ---
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(8).Hidden = $true
$ws.Rows.Item(9).Hidden = $true
$ws.Rows.Item(10).Hidden = $true
$ws.Rows.Item(11).Hidden = $true
$ws.Rows.Item(12).Hidden = $true
$ws.Rows.Item(13).Hidden = $true
$ws.Rows.Item(14).Hidden = $true
$ws.Rows.Item(15).Hidden = $true
$ws.Rows.Item(16).Hidden = $true
$ws.Rows.Item(17).Hidden = $true
$ws.Rows.Item(18).Hidden = $true
$ws.Rows.Item(19).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# --- Update the table's AutoFilter on the Class column: show CitationGeoKeys only ---
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(1, @("requirements_class_CitationGeoKeys"), 7)

# --- Move the selection/active cell to reflect the newly-scrolled view ---
$ws.Range("E4").Select()
